# Update attendance summary sheet:
#  - Column H (Absent) = 1 for rows 3,4,5,6,8,9,10,11,13,14,16,18
#  - Column G (Invalid) = 1 for row 3
#  - Columns D (Total Attendance Count) and E (Real) = 1 for rows 7, 12, 15, 17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H: Absent -> 1 for all rows except 7, 12, 15, 17
$rowsH = @(3, 4, 5, 6, 8, 9, 10, 11, 13, 14, 16, 18)
foreach ($r in $rowsH) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Column G: Invalid -> 1 for row 3
$ws.Cells.Item(3, 7).Value = 1

# Columns D and E -> 1 for rows 7, 12, 15, 17
$rowsDE = @(7, 12, 15, 17)
foreach ($r in $rowsDE) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
}
